{"js": "// Load all paragraphs in the body so we can locate the ones we need to touch.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) \"T\u00e1 dando onda\" is currently split across two runs with a\n//    <w:proofErr> pair in between (leftover grammar-check marks).\n//    Re-write the paragraph's content as a single clean run.\nconst target = paragraphs.items.find(p => p.text === \"T\u00e1 dando onda\");\nif (target) {\n    target.clear();\n    target.insertText(\"T\u00e1 dando onda\", Word.InsertLocation.start);\n}\n\n// 2) Add a new paragraph \"Mercen\u00e1rios\" right after \"American Pie\".\nconst americanPie = paragraphs.items.find(p => p.text === \"American Pie\");\nif (americanPie) {\n    americanPie.insertParagraph(\"Mercen\u00e1rios\", Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two paragraphs we need to touch by their text content so the\n# script does not depend on fixed paragraph indices.\n$target = $null\n$americanPie = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"T\") -and $t.Contains(\"dando onda\")) { $target = $p }\n    if ($t.StartsWith(\"American Pie\")) { $americanPie = $p }\n}\n\n# 1) \"Ta dando onda\" is split across two runs with a stray grammar-check\n#    <w:proofErr> pair in between. Insert a brand-new, clean paragraph with\n#    the same text right after it, then delete the old (marked-up) one.\nif ($target -ne $null) {\n    $r = $target.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $r.InsertAfter(\"T\u00e1 dando onda\")\n    $target.Range.Delete()\n}\n\n# 2) Add the new movie \"Mercenarios\" right after \"American Pie\".\nif ($americanPie -ne $null) {\n    $r2 = $americanPie.Range\n    $r2.Collapse(0)\n    $r2.InsertParagraphAfter()\n    $r2.Collapse(0)\n    $r2.InsertAfter(\"Mercen\u00e1rios\")\n}\n\n$d.Save()\n"}
